$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = "Dr. Gehan Adel, Dr. Amira Sobhy, Administrator, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat"
$ws.Cells.Item(3, 7).Value = "Administrator, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Veronia Rafat"
$ws.Cells.Item(4, 7).Value = "Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Majorelle Magdy"
$ws.Cells.Item(5, 7).Value = "Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Cells.Item(6, 7).Value = "Dr. Manar Montaser, Dr. Mohammad El-Tanany, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef"
$ws.Cells.Item(7, 7).Value = "Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Kerelos Zareef, Dr. Fatma Elhady, Dr. Lamiaa Ossama"
$ws.Cells.Item(9, 7).Value = "Dr. Safa Hany, Dr. Shimaa Ashraf"
$ws.Cells.Item(11, 7).Value = "Dr. Safa Hany, Dr. Aya Saeed, Dr. Amal Awwad"
$ws.Cells.Item(12, 7).Value = "Dr. Eman M. Abo-Sakaya, Dr. Amira Ibrahim, Dr. Marina Youhanna, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel"
$ws.Cells.Item(13, 7).Value = "Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh"
$ws.Cells.Item(15, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat"
$ws.Cells.Item(17, 7).Value = "Dr. Esraa Samy, Dr. Mohammad Safwat"
$ws.Cells.Item(19, 7).Value = "Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges"
$ws.Cells.Item(20, 7).Value = "Dr. Mariam Toma Gerges, Dr. Mohammad Safwat"
$ws.Cells.Item(28, 7).Value = "Dr. Maryam Ashraf, Dr. Aya Emad"
$ws.Cells.Item(30, 7).Value = "Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Shorok Mohammad"
